{"js": "const pairs = [\n  [\"2025-11-11 Tuesday\", \"2025-11-12 Wednesday\"],\n  [\"857\u00f72=\", \"328\u00f77=\"],\n  [\"468\u00f74=\", \"724\u00f75=\"],\n  [\"428\u00f78=\", \"677\u00f73=\"],\n  [\"516\u00f75=\", \"309\u00f78=\"],\n  [\"892\u00f73=\", \"418\u00f72=\"],\n  [\"768\u00f77=\", \"660\u00f72=\"],\n  [\"514\u00f73=\", \"300\u00f74=\"],\n  [\"461\u00f72=\", \"694\u00f75=\"],\n  [\"936\u00f73=\", \"978\u00f78=\"],\n  [\"540\u00f79=\", \"744\u00f79=\"],\n  [\"545\u00f77=\", \"562\u00f77=\"],\n  [\"212\u00f77=\", \"571\u00f77=\"],\n  [\"466\u00f72=\", \"613\u00f75=\"],\n  [\"879\u00f74=\", \"391\u00f78=\"],\n  [\"722\u00f73=\", \"797\u00f72=\"],\n  [\"152\u00f76=\", \"351\u00f79=\"],\n  [\"484\u00f72=\", \"829\u00f76=\"],\n  [\"320\u00f76=\", \"661\u00f74=\"],\n  [\"677\u00f78=\", \"888\u00f74=\"],\n  [\"230\u00f75=\", \"373\u00f78=\"],\n  [\"409\u00f75=\", \"471\u00f79=\"],\n  [\"498\u00f75=\", \"414\u00f75=\"],\n  [\"708\u00f79=\", \"795\u00f77=\"],\n  [\"904\u00f76=\", \"939\u00f72=\"],\n  [\"102\u00f77=\", \"487\u00f78=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-11-11 Tuesday\", \"2025-11-12 Wednesday\"),\n    @(\"857\u00f72=\", \"328\u00f77=\"),\n    @(\"468\u00f74=\", \"724\u00f75=\"),\n    @(\"428\u00f78=\", \"677\u00f73=\"),\n    @(\"516\u00f75=\", \"309\u00f78=\"),\n    @(\"892\u00f73=\", \"418\u00f72=\"),\n    @(\"768\u00f77=\", \"660\u00f72=\"),\n    @(\"514\u00f73=\", \"300\u00f74=\"),\n    @(\"461\u00f72=\", \"694\u00f75=\"),\n    @(\"936\u00f73=\", \"978\u00f78=\"),\n    @(\"540\u00f79=\", \"744\u00f79=\"),\n    @(\"545\u00f77=\", \"562\u00f77=\"),\n    @(\"212\u00f77=\", \"571\u00f77=\"),\n    @(\"466\u00f72=\", \"613\u00f75=\"),\n    @(\"879\u00f74=\", \"391\u00f78=\"),\n    @(\"722\u00f73=\", \"797\u00f72=\"),\n    @(\"152\u00f76=\", \"351\u00f79=\"),\n    @(\"484\u00f72=\", \"829\u00f76=\"),\n    @(\"320\u00f76=\", \"661\u00f74=\"),\n    @(\"677\u00f78=\", \"888\u00f74=\"),\n    @(\"230\u00f75=\", \"373\u00f78=\"),\n    @(\"409\u00f75=\", \"471\u00f79=\"),\n    @(\"498\u00f75=\", \"414\u00f75=\"),\n    @(\"708\u00f79=\", \"795\u00f77=\"),\n    @(\"904\u00f76=\", \"939\u00f72=\"),\n    @(\"102\u00f77=\", \"487\u00f78=\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
